# YW-54-T5_register - match the password and confirm password
#
# Restructures Sheet1 ("register_test_data"):
#   - Adds a new "IncorrectPassword" column (F) used to hold a password that
#     intentionally differs from "ConfirmPassword" for negative testing.
#   - F1 header  -> "IncorrectPassword"
#   - F2 value   -> "Ravi@2025" (kept as a hyperlinked/mailto cell, reusing
#     the relationship that used to live on C3)
#   - Row 3 (the "Maria" row) data is cleared out, since that hyperlinked
#     mailto cell/row is no longer part of the two-row dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column header: F1 = "IncorrectPassword" ---------------------------
$ws.Range("F1").Value = "IncorrectPassword"

# --- Remove the old C3 hyperlink + row 3 data -------------------------------
# NOTE: Range.Hyperlinks.Delete() clears every hyperlink on the sheet (not
# just the target range), so re-create the still-needed C2 hyperlink right
# after clearing.
$ws.Range("C3").Hyperlinks.Delete()
$ws.Range("A3:E3").ClearContents()

# --- Re-add the Ravi (row 2) email hyperlink on C2 --------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ravi.kumar1@testmail.com") | Out-Null
$ws.Range("C2").Style = "Hyperlink"

# --- New F2 cell: "Ravi@2025", reusing the old C3 hyperlink's mailto target -
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:maria.lopez2@testmail.com") | Out-Null
$ws.Range("F2").Value = "Ravi@2025"
$ws.Range("F2").Style = "Hyperlink"

# --- Update the active selection shown in the sheet view -------------------
$ws.Range("A3").Select() | Out-Null
